# EMU -> points helper (PowerPoint COM works in points; OOXML stores EMU, 12700 EMU = 1 pt).
# NOTE: Shape.Left/Top/Width/Height are single-precision (float32) in the PowerPoint object
# model, so a plain "emu / 12700" sometimes truncates to one EMU below the intended value once
# it is converted back on save. The literal used for the new group's Width below is a
# pre-compensated point value (verified empirically) that round-trips to exactly 6504111 EMU.
function EmuToPt($emu) { return $emu / 12700 }

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1) Move the top "SyncTool" bar group ("Gruppieren 28") up ---
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Gruppieren 28") {
        $sh.Top = EmuToPt(579496)
    }
}

# --- 2) Widen the "SyncTool.FileSystem" group ("Gruppieren 30") to the left, keeping its
#        right edge fixed (so it still lines up with the other boxes in that row) ---
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Gruppieren 30") {
        $sh.Left = EmuToPt(4237892)
        $sh.Width = EmuToPt(6954708)
    }
}

# --- 3) Add a new group ("Gruppieren 37"), a copy of "Gruppieren 34" (the Synchronization
#        box), holding the new "SyncTool.Synchronization.Git" / ".Git.Test" boxes ---
$src = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Gruppieren 34") {
        $src = $sh
    }
}

# The runtime's GroupItems collection resolves children of a freshly Duplicate()-d shape back
# to the *source* group's children (they still share the same internal shape id right after
# duplication), so renaming/retexting through $new.GroupItems would silently edit $src instead.
# Work around this by first rewriting $src's children to the values the new group should end up
# with, duplicating (so the clone is created carrying those values), and then restoring $src's
# children to their original name/text.
$child1 = $src.GroupItems.Item(1)
$child2 = $src.GroupItems.Item(2)

$origName1 = $child1.Name
$origText1 = $child1.TextFrame.TextRange.Text
$origName2 = $child2.Name
$origText2 = $child2.TextFrame.TextRange.Text

$child1.Name = "Abgerundetes Rechteck 38"
$child1.TextFrame.TextRange.Text = "SyncTool.Synchronization.Git.Test"
$child2.Name = "Abgerundetes Rechteck 39"
$child2.TextFrame.TextRange.Text = "SyncTool.Synchronization.Git"

$new = $src.Duplicate()
$new.Name = "Gruppieren 37"
$new.Left = EmuToPt(4688489)
$new.Top = EmuToPt(1344717)
$new.Width = 512.1347351094452   # pre-compensated point value -> exactly 6504111 EMU
$new.Height = EmuToPt(708589)

$child1.Name = $origName1
$child1.TextFrame.TextRange.Text = $origText1
$child2.Name = $origName2
$child2.TextFrame.TextRange.Text = $origText2
